$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, after the weekly re-ordering/update of the price records.
# Columns: D=Fecha, J=Volumen, K=Precio minimo, L=Precio maximo, M=Precio promedio ponderado, P=Precio $/Kg
$rows = @(
    @{ Row = 2; D = 45084; J = 90;  K = 22000; L = 23000; M = 22556; P = 1504 },
    @{ Row = 3; D = 44749; J = 90;  K = 17000; L = 18000; M = 17556; P = 1170 },
    @{ Row = 4; D = 45063; J = 40;  K = 21000; L = 22000; M = 21500; P = 1433 },
    @{ Row = 5; D = 45091; J = 40;  K = 20000; L = 22000; M = 21000; P = 1400 },
    @{ Row = 6; D = 44750; J = 140; K = 19000; L = 20000; M = 19571; P = 1305 },
    @{ Row = 7; D = 44839; J = 50;  K = 15000; L = 16000; M = 15600; P = 1040 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.D    # D - Fecha
    $ws.Cells.Item($row, 10).Value = $r.J   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $r.P   # P - Precio $/Kg
}
